# "ajout de la redirection" - add two new invitee rows (20 and 21)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "isaac"
$ws.Range("B20").Value = "Omar"
$ws.Range("C20").Value = 74160680
$ws.Range("D20").Value = "gabon"
$ws.Range("E20").Value = "assekoazareel222@gmail.com"

$ws.Range("A21").Value = "gsvxshcvsh"
$ws.Range("B21").Value = "hsxsxs"
$ws.Range("C21").Value = 74951291
$ws.Range("D21").Value = "gabon"
$ws.Range("E21").Value = "assekoazareel222@gmail.com"
